$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.564.40'
$ws.Range("E2").Value = '  -0.74%  '
$ws.Range("D3").Value = '2.072.86'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.21'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.92'
$ws.Range("E8").Value = '  -1.85%  '
$ws.Range("E9").Value = '  -1.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0776'
$ws.Range("E10").Value = '  -1.15%  '
$ws.Range("E11").Value = '  +1.64%  '
$ws.Range("D12").Value = '2.373.99'
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.76'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.36'
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.765'
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.37'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").Value = '2.065.26'
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("D18").Value = '37.529.46'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.17'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.90'
$ws.Range("E20").Value = '  -2.27%  '
$ws.Range("E21").Value = '  -2.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.67'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.40'
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("E25").Value = '  -2.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.90'
$ws.Range("E26").Value = '  +7.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.76'
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.132'
$ws.Range("E28").Value = '  -4.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.29'
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.36'
$ws.Range("E30").Value = '  -3.87%  '
$ws.Range("E31").Value = '  +0.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.57'
$ws.Range("E32").Value = '  -3.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0627'
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.62'
$ws.Range("E34").Value = '  -2.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.53'
$ws.Range("E35").Value = '  +1.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.82'
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.31'
$ws.Range("E37").Value = '  -3.17%  '
$ws.Range("E38").Value = '  +0.19%  '
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0227'
$ws.Range("E40").Value = '  +4.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.04'
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0963'
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.21'
$ws.Range("E43").Value = '  +4.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.90'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").Value = '1.479.87'
$ws.Range("E45").Value = '  +2.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.75'
$ws.Range("E46").Value = '  -0.80%  '
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.01'
$ws.Range("E48").Value = '  -3.82%  '
$ws.Range("E49").Value = '  -1.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.97'
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").Value = '2.258.91'
$ws.Range("E51").Value = '  -0.75%  '
